$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the existing row 2 (the old "Yao / Water"
# triplicate block). This pushes the old rows 2-8 down to rows 4-10 and
# gives us 5 rows (2-6) to hold the new "Dee" investigation dialogue.
$ws.Rows.Item(2).Resize(2).Insert()

# --- Row 2: Dee - "Upon re-examining..." (Suspicious / appearAt) ---
$ws.Range("A2").Value = "Dee"
$ws.Range("B2").Value = "Upon re-examining the area near the main door, we found tiny traces of blood on the threshold."
$ws.Range("C2").Value = "Dee-Determined"
$ws.Range("D2").Value = "DialogueVocal"
$ws.Range("E2").Value = "DoorInvestigate"
$ws.Range("F2").Value = "Suspicious"
$ws.Range("J2").Value = "appearAt"
$ws.Range("K2").Value = 500
$ws.Range("L2").Value = "Dee-Thinking"

# --- Row 3: Dee - "There are only a few drops..." ---
$ws.Range("A3").Value = "Dee"
$ws.Range("B3").Value = "There are only a few drops, very faint, but the color and condition match the pool on the ground——"
$ws.Range("C3").Value = "Dee-Thinking2"
$ws.Range("D3").Value = "DialogueVocal"
$ws.Range("E3").Value = "DoorInvestigate"

# --- Row 4 (was old row 2 "Yao/Water"): Dee - "Likely left at the same time." ---
$ws.Range("A4").Value = "Dee"
$ws.Range("B4").Value = "Likely left at the same time."
$ws.Range("C4").Value = "Dee-Thinking2"
$ws.Range("D4").Value = "DialogueVocal"
$ws.Range("E4").Value = "DoorInvestigate"

# --- Row 5 (was old row 3 "Yao/Water"): Dee - "These must be the Lord's bloodstains." ---
$ws.Range("A5").Value = "Dee"
$ws.Range("B5").Value = "These must be the Lord" + [char]0x2019 + "s bloodstains."
$ws.Range("C5").Value = "Dee-Determined"
$ws.Range("D5").Value = "DialogueVocal"
$ws.Range("E5").Value = "DoorInvestigate"

# --- Row 6 (was old row 4 "Yao/Water"): Dee - "The shape......it's rather intriguing." ---
$ws.Range("A6").Value = "Dee"
$ws.Range("B6").Value = "The shape......it" + [char]0x2019 + "s rather intriguing."
$ws.Range("C6").Value = "Dee-Thinking2"
$ws.Range("D6").Value = "DialogueVocal"
$ws.Range("E6").Value = "DoorInvestigate"

# Rows 2 and 3 hold the longer wrapped lines, so they need the taller row height.
$ws.Range("A2").EntireRow.RowHeight = 34
$ws.Range("A3").EntireRow.RowHeight = 34

# --- Row 7 (was old row 5 "Investigate2/Water/Water/DialogueVocal/DoorInvestigate") ---
# Keep A7/B7/C7, but clear the now-unused D7/E7 (DialogueVocal/DoorInvestigate).
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()

# --- Row 8 (was old row 6 "Hand/Hand/DialogueVocal/DoorInvestigate") ---
# Keep B8/C8, clear D8/E8, and add the J/K/L highlight cells (style only, no value).
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("J8").Value = $ws.Range("J7").Value
$ws.Range("J8").ClearContents()
$ws.Range("J8").Style = $ws.Range("J9").Style
$ws.Range("K8").Style = $ws.Range("K9").Style
$ws.Range("L8").Style = $ws.Range("L9").Style

# Rows 9 (Blood/Blood/disappear) and 10 (End Investigation/StoryScript14) already
# carry the correct values/styles after the row shift - no further edits needed.

# Restore the active selection to match the saved workbook.
$ws.Range("B8").Select()
